# Update cryptocurrency price/volume data, preserving the sheet's
# convention of storing Price/Volume figures as plain text (no
# locale re-interpretation as numbers/dates).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, [string]$addr, [string]$val)
    $rng = $ws.Range($addr)
    # Force text interpretation so values like "0.160" or "15.10"
    # keep their trailing zeros instead of being auto-parsed as numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    # Drop the now-unneeded "Text" number format / quote-prefix style so
    # the cell's style index matches the rest of the (unstyled) data cells.
    $rng.ClearFormats()
}

# --- Row 47 / 48: EnergySwap and InjectiveProtocol swapped ranking order ---
Set-TextCell $ws "B47" "InjectiveProtocol"
Set-TextCell $ws "C47" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws "D47" "25.06"
Set-TextCell $ws "E47" "  -3.29%  "

Set-TextCell $ws "B48" "EnergySwap"
Set-TextCell $ws "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D48" "24.29"
Set-TextCell $ws "E48" "  +3.37%  "

# --- Price (D) / Volume(1h) (E) refresh for all other rows ---
Set-TextCell $ws "D2" "66.646.08"
Set-TextCell $ws "E2" "  +0.73%  "
Set-TextCell $ws "D3" "3.596.99"
Set-TextCell $ws "E3" "  +0.96%  "
Set-TextCell $ws "E4" "  -0.03%  "
Set-TextCell $ws "D5" "609.75"
Set-TextCell $ws "E5" "  +0.40%  "
Set-TextCell $ws "D6" "148.22"
Set-TextCell $ws "E6" "  +2.34%  "
Set-TextCell $ws "E7" "  +0.08%  "
Set-TextCell $ws "E8" "  +0.45%  "
Set-TextCell $ws "D9" "8.05"
Set-TextCell $ws "E9" "  +1.25%  "
Set-TextCell $ws "E10" "  -0.16%  "
Set-TextCell $ws "D12" "4.206.17"
Set-TextCell $ws "E12" "  +0.91%  "
Set-TextCell $ws "E13" "  +0.81%  "
Set-TextCell $ws "D14" "29.96"
Set-TextCell $ws "E14" "  -0.75%  "
Set-TextCell $ws "D15" "3.642.61"
Set-TextCell $ws "E15" "  +2.20%  "
Set-TextCell $ws "D16" "66.719.79"
Set-TextCell $ws "E16" "  +0.67%  "
Set-TextCell $ws "E17" "  +0.81%  "
Set-TextCell $ws "D18" "11.51"
Set-TextCell $ws "E18" "  +0.40%  "
Set-TextCell $ws "D19" "6.33"
Set-TextCell $ws "E19" "  +1.94%  "
Set-TextCell $ws "D20" "15.10"
Set-TextCell $ws "E20" "  +1.22%  "
Set-TextCell $ws "D21" "427.87"
Set-TextCell $ws "E21" "  -0.79%  "
Set-TextCell $ws "E22" "  +1.61%  "
Set-TextCell $ws "D23" "78.95"
Set-TextCell $ws "E23" "  +0.15%  "
Set-TextCell $ws "D24" "3.737.20"
Set-TextCell $ws "E24" "  +0.87%  "
Set-TextCell $ws "E25" "  +0.00%  "
Set-TextCell $ws "E26" "  +3.51%  "
Set-TextCell $ws "D27" "8.31"
Set-TextCell $ws "E27" "  +3.52%  "
Set-TextCell $ws "D28" "9.32"
Set-TextCell $ws "E28" "  +2.24%  "
Set-TextCell $ws "E29" "  -0.12%  "
Set-TextCell $ws "E30" "  -0.01%  "
Set-TextCell $ws "D31" "0.160"
Set-TextCell $ws "E31" "  +1.35%  "
Set-TextCell $ws "D32" "3.594.53"
Set-TextCell $ws "E32" "  +0.99%  "
Set-TextCell $ws "D33" "1.47"
Set-TextCell $ws "E33" "  -1.13%  "
Set-TextCell $ws "E34" "  -0.08%  "
Set-TextCell $ws "E35" "  -0.89%  "
Set-TextCell $ws "E36" "  +0.04%  "
Set-TextCell $ws "E37" "  +0.22%  "
Set-TextCell $ws "E38" "  -2.73%  "
Set-TextCell $ws "D39" "177.58"
Set-TextCell $ws "E39" "  +4.40%  "
Set-TextCell $ws "E40" "  +0.30%  "
Set-TextCell $ws "E41" "  +0.35%  "
Set-TextCell $ws "E42" "  +0.19%  "
Set-TextCell $ws "D43" "1.91"
Set-TextCell $ws "E43" "  -1.36%  "
Set-TextCell $ws "D44" "2.57"
Set-TextCell $ws "E44" "  +8.34%  "
Set-TextCell $ws "E45" "  +0.00%  "
Set-TextCell $ws "E46" "  -1.44%  "
Set-TextCell $ws "E49" "  +0.55%  "
Set-TextCell $ws "D50" "0.953"
Set-TextCell $ws "E50" "  +0.20%  "
Set-TextCell $ws "D51" "0.236"
Set-TextCell $ws "E51" "  -1.40%  "
